$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The runs/balls/fours columns (C:E) hold numeric-looking values that are
# actually stored as text in the source data. Force the text number
# format on just the cells being touched, so re-entering the values
# doesn't get auto-coerced to numbers by Excel.
$ws.Range("C2:E6").NumberFormat = "@"
$ws.Range("C8:E8").NumberFormat = "@"
$ws.Range("C10:E12").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = "6"
$ws.Range("D2").Value = "6"
$ws.Range("E2").Value = "1"

# Row 3
$ws.Range("C3").Value = "0"
$ws.Range("D3").Value = "2"

# Row 4
$ws.Range("C4").Value = "12"
$ws.Range("D4").Value = "13"
$ws.Range("E4").Value = "0"

# Row 5
$ws.Range("C5").Value = "11"
$ws.Range("D5").Value = "7"
$ws.Range("E5").Value = "1"

# Row 6
$ws.Range("C6").Value = "10"
$ws.Range("D6").Value = "5"
$ws.Range("E6").Value = "2"

# Row 8
$ws.Range("C8").Value = "32"
$ws.Range("D8").Value = "24"
$ws.Range("E8").Value = "3"

# Row 10
$ws.Range("C10").Value = "11"
$ws.Range("D10").Value = "18"
$ws.Range("E10").Value = "0"

# Row 11
$ws.Range("C11").Value = "1"
$ws.Range("D11").Value = "4"
$ws.Range("E11").Value = "0"

# Row 12
$ws.Range("C12").Value = "7"
$ws.Range("D12").Value = "12"
